# Adds 10 new hackathon song submissions (rows 17-26), a thin separator
# row, and re-numbers the remaining "placeholder" counting rows below -
# including the accidental duplicate row that exists in the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen column B to fit the longer song names being added ------------
# (31.1640625 is the exact authored width; 30.25 is the closest input this
# host's column-width rounding grid can land on)
$ws.Columns("B").ColumnWidth = 30.25

# --- new data rows 17-26 (overwrite the old placeholder rows in place) --
$newRows = @(
    @(17, "I´m the plug",                    "Drake",          48557410, "0.88725203", 2, "Mike "),
    @(18, "Stepped on my j´z",                "Nelly",          4359488,  "0.7798307",  2, "Mike "),
    @(19, "Magic Bus",                        "The Whoo",       4735362,  "0.88725203", 1, "Mike "),
    @(20, "The electric co",                  "U2",             3128392,  "0.80582285", 2, "Mike "),
    @(21, "Saviour",                          "Rise against",   49697071, "0.85952747", 2, "Mike "),
    @(22, "Red Red Wine",                     "Neil Diamond",   37858536, "1.1771723",  1, "Mike "),
    @(23, "Aint no grave",                    "Johnny Cash",    8216660,  "1.1954225",  1, "Mike "),
    @(24, "At last",                          "Etta Janes",     5510572,  "1.0348469",  1, "Mike "),
    @(25, "Amazing Grace",                    "Nana Mouskouri", 596677,   "0.83238655", 1, "Mike "),
    @(26, "Happy Happy Happy Wedding Day",    "Bing Crosby",    4667290,  "1.1400715",  1, "Mike ")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $r
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
}

# --- thin grey/black separator row, like the one used for column F ------
$ws.Rows(27).Insert()
$ws.Rows(27).RowHeight = 7
$ws.Range("A27:H27").Interior.Pattern = 17

# --- duplicate row 28 (value 27) exactly as found upstream --------------
$ws.Rows(30).Insert()
$ws.Cells.Item(30, 1).Value = 27

# --- selection, as recorded in the workbook's last-saved view -----------
$ws.Range("G26").Select()

Write-Output "edit applied"
